# Mark a batch of SYS-FLC requirements as "Complete" in the MASTER SPREADSHEET
# sheet, mirroring the author's manual edit/testing pass (FLC-012 .. FLC-083).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MASTER SPREADSHEET")

# Row 35..50 block (uses the same visual "Complete" look as the existing
# E104/E105 cells - light green fill, status text "Complete").
$rowsGroupA = @(35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 49, 50, 106)

# Row 55..98 block (uses the same visual "Complete" look as the existing
# E54/E60/E66... cells - light green fill, status text "Complete").
$rowsGroupB = @(55, 56, 62, 63, 68, 74, 80, 86, 92, 98)

foreach ($r in $rowsGroupA) {
    $target = $ws.Cells.Item($r, 5)
    $ws.Range("E104").Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
    $target.Value = "Complete"
}

foreach ($r in $rowsGroupB) {
    $target = $ws.Cells.Item($r, 5)
    $ws.Range("E54").Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
    $target.Value = "Complete"
}

# Minor formatting touch-up on E34 (still "Incomplete") picked up alongside
# the adjacent edits - match the thin-border look used elsewhere (e.g. E11).
$e34 = $ws.Range("E34")
$ws.Range("E11").Copy() | Out-Null
$e34.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the saved view state (scroll position / zoom / selection) to match
# where the author left off after testing these requirements.
$ws.Activate()
$ws.Range("D109").Select() | Out-Null
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1
